$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet: bump version/status/date/contact for the new release ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value  = "0.4.0-snapshot-1"               # Version
$meta.Range("B6").Value  = "draft"                          # Status
$meta.Range("B8").Value  = "2024-05-23T12:16:26+00:00"      # Date
$meta.Range("B10").Value = "ANS (https://esante.gouv.fr)"   # Contact

# --- "Elements" sheet: the two mapping columns (AK = col 37, AL = col 38)
#     swap places - "Mapping: Spécification métier..." now comes first,
#     followed by "Mapping: RIM Mapping". Swap the cell contents (rather
#     than cut/insert, which would disturb every other column's width)
#     and then swap the two columns' widths to match.
$els = $wb.Worksheets.Item("Elements")

$xlUp = -4162
$lastRow = $els.Cells.Item($els.Rows.Count, 37).End($xlUp).Row
for ($r = 1; $r -le $lastRow; $r++) {
    $leftCell  = $els.Cells.Item($r, 37)
    $rightCell = $els.Cells.Item($r, 38)
    $leftVal   = $leftCell.Value2
    $rightVal  = $rightCell.Value2
    if ($leftVal -ne $rightVal) {
        $leftCell.Value  = $rightVal
        $rightCell.Value = $leftVal
    }
}

$els.Columns.Item(37).ColumnWidth = 81.15
$els.Columns.Item(38).ColumnWidth = 24.15
